$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Row 2
Set-TextCell "D2" "300.43"
Set-TextCell "E2" "0.50%"

# Row 3
Set-TextCell "D3" "32.19"
Set-TextCell "E3" "1.79%"

# Row 4
Set-TextCell "D4" "4.956"
Set-TextCell "E4" "-3.75%"

# Row 5
Set-TextCell "D5" "0.07876"

# Row 6
Set-TextCell "D6" "2.091"
Set-TextCell "E6" "-14.98%"

# Row 7
Set-TextCell "D7" "7.795"
Set-TextCell "E7" "-0.07%"

# Row 8
Set-TextCell "D8" "3.837"
Set-TextCell "E8" "-2.00%"

# Row 9
Set-TextCell "D9" "0.9262"
Set-TextCell "E9" "-0.28%"

# Row 10
Set-TextCell "D10" "0.1741"
Set-TextCell "E10" "-1.52%"

# Row 11
Set-TextCell "D11" "0.07952"
Set-TextCell "E11" "7.09%"

# Row 12
Set-TextCell "D12" "0.08643"
Set-TextCell "E12" "-2.41%"

# Row 13
Set-TextCell "D13" "0.03105"
Set-TextCell "E13" "3.36%"

# Row 14
Set-TextCell "D14" "0.1002"
Set-TextCell "E14" "0.17%"

# Row 15
Set-TextCell "D15" "0.001528"
Set-TextCell "E15" "0.84%"

# Row 16
Set-TextCell "D16" "0.005794"
Set-TextCell "E16" "-4.16%"

# Row 17
Set-TextCell "E17" "2,099.89%"

# Row 18
Set-TextCell "D18" "3.460"
Set-TextCell "E18" "-1.93%"

# Row 19
Set-TextCell "D19" "2.243"
Set-TextCell "E19" "-2.03%"

# Row 20
Set-TextCell "E20" "0.46%"

# Row 21
Set-TextCell "D21" "0.1310"
Set-TextCell "E21" "-2.18%"

# Row 22
Set-TextCell "D22" "4.303"
Set-TextCell "E22" "3.64%"

# Row 23
Set-TextCell "E23" "6.80%"

# Row 24
Set-TextCell "D24" "0.04613"
Set-TextCell "E24" "-0.18%"

# Row 25
Set-TextCell "D25" "0.001236"
Set-TextCell "E25" "-0.44%"

# Row 26
Set-TextCell "D26" "0.004428"
Set-TextCell "E26" "-2.20%"

# Row 39
Set-TextCell "D39" "0.01708"
Set-TextCell "E39" "-2.35%"

# Row 40
Set-TextCell "D40" "0.04752"
Set-TextCell "E40" "3.49%"

# Row 41
Set-TextCell "D41" "0.007421"
Set-TextCell "E41" "7.37%"

# Row 42
Set-TextCell "D42" "0.1355"
Set-TextCell "E42" "-1.32%"

# Row 43
Set-TextCell "D43" "0.002354"
Set-TextCell "E43" "7.49%"

# Row 44
Set-TextCell "D44" "0.01133"
Set-TextCell "E44" "10.10%"

# Row 45
Set-TextCell "D45" "0.00006024"
Set-TextCell "E45" "-3.88%"

# Row 46
Set-TextCell "E46" "0.15%"

# Row 47
Set-TextCell "E47" "-59.58%"

# Row 48
Set-TextCell "D48" "0.8204"
Set-TextCell "E48" "9.61%"

# Row 49
Set-TextCell "E49" "0.15%"

# Row 50
Set-TextCell "E50" "0.15%"
